$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
  'Kaleida - Think (Lyrics) John Wick soundtrack',
  'Hands Like Houses - Torn',
  'E.Satie - Gnossienne N.1 (Piano)',
  'La Tormenta De Arena - Dorian (letra)',
  'Yanni - Can''t Wait (Sensuous Chill)',
  'Zivert - Life (English Version)',
  'Mariah Carey - My All [Lyrics]',
  'One Direction - Story of My Life',
  'Beth Thornton - Something You Don''t Know',
  '2CELLOS - Shape Of My Heart [Live at Arena di Verona]',
  'Enrique Iglesias - EL BAÑO (Letra) ft. Bad Bunny',
  '2CELLOS - Love Story',
  'Lana Del Rey - Video Games',
  '2CELLOS - Fragile [LIVE at Arena Pula]',
  'Justin Bieber & benny blanco - Lonely (Official Acoustic Video)',
  '서태지와 아이들   이 밤이 깊어가지만 (가사 첨부)',
  'Cody Francis - Rose In The Garden',
  'I Got Summer On My Mind (Still Dre Remix)',
  'Lady Gaga, Bruno Mars - Die With A Smile',
  'Maxim Fadeev - Googoosha',
  'Elley Duhe - Middle Of The Night',
  'Broken Angel (Albert Vishi ft. Taulant Sllamniku Cover)',
  'Laura Pausini - It’s Not Goodbye',
  'Juice Wrld - Lucid Dreams ',
  'Ramz - Barking',
  'falling in love with someone you can''t have (a playlist)',
  'Feeling Good 🌱 A playlist to lift your mood',
  'Best classical music: Beethoven, Mozart, Schubert,Bach...🎶',
  'Selena Gomez - Bad Liar',
  'Florida Georgia Line - Simple (Lyrics)',
  'Tones and I: Dance Monkey (US TV Debut)',
  'the luka state - bring this all together',
  'Sunset Sons - I Can`t Wait (Official Audio)',
  'Sunset Sons - Somewhere Maybe (Official Audio)',
  'Sunset Sons - Know My Name (Official Audio)',
  'Sunset Sons - The River',
  'Sunset Sons - Loa (Official Audio)',
  'Sunset Song - On The Road (Lyrics)',
  'Sunset Sons - Remember',
  'Craig David - Walking Away [Lyrics] 🎵',
  'Escape (Rosaline OST)',
  'Zara Larsson - Lush Life',
  'JEON SOMI DUMB DUMB Lyrics (전소미 DUMB DUMB 가사)',
  'Shakira - Can`t Remember to Forget You (Lyrics) ft. Rihanna',
  'Havana feat. Yaar &amp; Kaiia - Last Night (Lyrics)',
  'Craig David - Rise & Fall ft. Sting (Official Video)',
  'Maroon 5 - Girls Like You ft. Cardi B (Official Music Video)',
  'Havana feat. Yaar & Kaiia - Big Love (Official Video)',
  'Edward Maya, Vika Jigulina - Stereo love (Radio Edit) (Lyrics)',
  'Gym Class Heroes: Stereo Hearts ft. Adam Levine',
  'Shawn Mendes - It''ll Be Okay',
  'somewhere only we know (Gustixa &amp; Rhianne)',
  'Duncan Laurence feat. FLETCHER – Arcade',
  'Bad Liar (cover)',
  'Lana Del Rey - Summertime Sadness (Official Music Video)',
  'Linkin Park - In The End (Mellen Gi &amp; Tommee Profitt Remix)',
  'Enya - Only Time (Official 4K Music Video)',
  'Today is a Good Day',
  'Heroes Tonight',
  'Shakira - Hips Dont Lie',
  'Shakira - La La La World Cup 2014',
  'Something Just Like This',
  'Shakira - Chantaje (letra)',
  'Zara Larsson – Dont Worry Bout Me',
  'Selena Gomez - Buscando Amor',
  'Shawn Mendes - In My Blood',
  'Charlie Puth - We Dont Talk Anymore',
  'Selena Gomez - Adiós',
  'Charlie Puth - Attention',
  'Bruno Mars - Grenade',
  'Let Me Love You',
  'Sweet But Psycho',
  'Who`s laughing now - Ava Max',
  'Girls Like You',
  'Camila Cabello - Havana',
  'Maroon 5 - Memories',
  'post malone - rockstar (feat. 21 savage)',
  'Drake - Gods plan',
  'Dua Lipa - Levitating',
  'Justin Bieber - Baby',
  'State of Grace',
  'Taylor Swift - Back To December',
  'As it was',
  'Taylor Swift - Begin Again',
  'Sofia',
  'Numb (lyrics|rock)',
  'For The Rest Of My Life',
  'Insha Allah',
  'Hunger Games | Atlas',
  'Hard To Say Im Sorry',
  'Sasha Sloan - Lie',
  'Solo Para Ti',
  'Zara Larsson - This Ones For You',
  'La Cintura',
  'Diamond Heart',
  'Halsey - Sorry',
  'Relax, Take it Easy',
  'Thank you Allah',
  'Let me down slowly',
  'Lonely (acoustic)',
  'Selfish love',
  'twenty one pilots: Heathens',
  'Dont Let Me Down',
  'Someone you loved',
  'Fed up with us',
  'Skyfall',
  'Story of my life',
  'Reamonn - Tonight',
  'Hymn For The Weekend',
  'Hell Or High Water',
  'Everything I Need',
  'Love the way you lie',
  'Moving Mountains',
  'Numb (cover)',
  'Until I found you',
  'Heart is on fire',
  'Holes',
  'Survivors',
  'Irakliy - Ya s toboy(cover)',
  'Vetrom stat (cover)',
  'Reamonn - Supergirl',
  'Another Love',
  'All the little lights',
  'Arcade',
  'Burito - Po volnam',
  'To Be Free',
  'Castle of Glass',
  'Ava Max - Alone',
  'Let Her Go (ft Ed Sheeran)',
  'In the end (rmx)',
  'No Time To Die',
  'The Cup Of Life - Ricky Martin',
  'Баста - Выпускной',
  'Wonderful Life',
  'Shape of My Heart',
  'Wonderful Life',
  'Shape of My Heart',
  'Tired Of Being Sorry'
)

$links = @(
  'https://www.youtube.com/watch?v=FVtFcbBfNYw',
  'https://www.youtube.com/watch?v=M58IJO7N32s',
  'https://www.youtube.com/watch?v=X3JLOenXGUc',
  'https://www.youtube.com/watch?v=28W-KrHjmK8',
  'https://www.youtube.com/watch?v=9kardLhsFrk',
  'https://www.youtube.com/watch?v=mTecGII7cFA',
  'https://www.youtube.com/watch?v=o4che1p-M4M',
  'https://www.youtube.com/watch?v=W-TE_Ys4iwM',
  'https://www.youtube.com/watch?v=RftohIbwlqg',
  'https://www.youtube.com/watch?v=jx1-NP9_YIA',
  'https://www.youtube.com/watch?v=8BbtBnnnvCM',
  'https://www.youtube.com/watch?v=UdHopftQD3A',
  'https://www.youtube.com/watch?v=cE6wxDqdOV0',
  'https://www.youtube.com/watch?v=q_ymIjOyzRQ',
  'https://www.youtube.com/watch?v=Cu5hhxP_prE',
  'https://www.youtube.com/watch?v=__SXVP2GmvM',
  'https://www.youtube.com/watch?v=JO4-j1LfoQQ',
  'https://www.youtube.com/watch?v=89LOsf8pDhY',
  'https://www.youtube.com/watch?v=zgaCZOQCpp8',
  'https://www.youtube.com/watch?v=gqOoJXttEec',
  'https://www.youtube.com/watch?v=KLTMCPzRO64',
  'https://www.youtube.com/watch?v=5miHGQVFJm0',
  'https://www.youtube.com/watch?v=onYQkI8S1UY',
  'https://www.youtube.com/watch?v=_fh64GbFSw4',
  'https://www.youtube.com/watch?v=Q0QKUU95bVc',
  'https://www.youtube.com/watch?v=_K57AlI62V4',
  'https://www.youtube.com/watch?v=VYtBO_cDJCU',
  'https://www.youtube.com/watch?v=DxnDcH2NS5c',
  'https://www.youtube.com/watch?v=NZKXkD6EgBk',
  'https://www.youtube.com/watch?v=TuTDc9d_9yI',
  'https://www.youtube.com/watch?v=4iQxG8ZjYO8',
  'https://www.youtube.com/watch?v=OcJ5EgxsWBg',
  'https://www.youtube.com/watch?v=JuiegvRQ8dI',
  'https://www.youtube.com/watch?v=SHapfmLyBp0',
  'https://www.youtube.com/watch?v=orMwK0veDVQ',
  'https://www.youtube.com/watch?v=MCyEm1fViZQ',
  'https://www.youtube.com/watch?v=9tXWQy7mMsM',
  'https://www.youtube.com/watch?v=NsKZ-5EDqPA',
  'https://www.youtube.com/watch?v=PH_P12XqY9Y',
  'https://www.youtube.com/watch?v=8AwamgSDpdA',
  'https://www.youtube.com/watch?v=M9b_z-LKE14',
  'https://www.youtube.com/watch?v=tD4HCZe-tew',
  'https://www.youtube.com/watch?v=TfAzTYzBvTo',
  'https://www.youtube.com/watch?v=i_XM3u1_jZQ',
  'https://www.youtube.com/watch?v=i-Yuf5-zTec',
  'https://www.youtube.com/watch?v=pU2ukeS2JTE',
  'https://www.youtube.com/watch?v=aJOTlE1K90k',
  'https://www.youtube.com/watch?v=aVFNJBqj5vU',
  'https://www.youtube.com/watch?v=y9Kqb2z9Lzs',
  'https://www.youtube.com/watch?v=T3E9Wjbq44E',
  'https://youtu.be/KrgJp7Z1Hv8?si=MOyY5rZzP-7kcfhM',
  'https://www.youtube.com/watch?v=92izkAK5OA0',
  'https://www.youtube.com/watch?v=308v08mFWWc',
  'https://youtu.be/5jfz3q9Z0RY?si=OHvyb7AMtM_wtAXc',
  'https://www.youtube.com/watch?v=TdrL3QxjyVw',
  'https://www.youtube.com/watch?v=WNeLUngb-Xg',
  'https://www.youtube.com/watch?v=7wfYIMyS_dI',
  'https://youtu.be/9L4EjJqrz0c?si=x97RAvAA9IELRZPW',
  'https://www.youtube.com/watch?v=074rfF4RJZc',
  'https://youtu.be/p3pEe6aAJ4k?si=bzrAEs7c-zSwqBUo',
  'https://youtu.be/2igups6VdcA?si=N5uu5genirJuWXWC',
  'https://youtu.be/FM7MFYoylVs?si=TrbAGj-JAUeEJ4bd',
  'https://youtu.be/J76eQJP3UIQ?si=juYKqG_UCEta8y19',
  'https://youtu.be/u_tzZd9kIWg?si=y-s2yCVh4U2JLsJJ',
  'https://youtu.be/2P6EExu3H5s?si=f2hv9y52VqxnVOmL',
  'https://youtu.be/36tggrpRoTI?si=CiCfVdO8Oepjt4Rs',
  'https://youtu.be/bpFVJJBgtXY?si=L2NuwOWGhmdKacwg',
  'https://youtu.be/9H_368c2Hzw?si=UOBGyTGbUe_fISFW',
  'https://youtu.be/Oz5JDtkf1as',
  'https://youtu.be/4YrzJ9RZ9qY',
  'https://youtu.be/SMs0GnYze34?si=T-UORWGqJCoitcOM',
  'https://youtu.be/2KBFD0aoZy8',
  'https://youtu.be/4JYSgIiSZSA?si=3v9kDuzvYJvWaOsO',
  'https://youtu.be/aJOTlE1K90k',
  'https://youtu.be/HCjNJDNzw8Y?si=QjZAi7GPIc4ParOQ',
  'https://www.youtube.com/watch?v=SlPhMPnQ58k&pp=ygUPbWFyb29uIG1lbW9yaWVz',
  'https://www.youtube.com/watch?v=9lQP9-F8kIQ',
  'https://www.youtube.com/watch?v=ScfgOVJiu_I',
  'https://www.youtube.com/watch?v=j2c3tR_qfiQ',
  'https://www.youtube.com/watch?v=khOFw2f4bQY',
  'https://www.youtube.com/watch?v=gr4cqcqnAN0',
  'https://youtu.be/QUwxKWT6m7U?si=LNPBWKl0DqXIfOP2',
  'https://youtu.be/Qfm6nfz1QNQ?si=3mMjYFpALij7GELl',
  'https://youtu.be/cMPEd8m79Hw?si=9zE5-51p0xGyEgSO',
  'https://youtu.be/ftI_Lp7LAuU?si=aOFT5Ral2-A_2PxG',
  'https://youtu.be/8P0vKLHbtMg?si=HhXMHjE8vD2yeC_B',
  'https://youtu.be/PHbZ9SXHJwA?si=_7a2Gaka2oPEWrCQ',
  'https://youtu.be/8xXJyFNfiy8?si=XkqgGm4hEyZoqJe1',
  'https://youtu.be/Lh3TokLzzmw?si=I5CcdBNIEuwDZvVT',
  'https://youtu.be/XCmOdVia9DE?si=60M6i15UUakuL7DH',
  'https://youtu.be/AzjTJpzfB8U?si=PHYxAGETm1P1opd0',
  'https://youtu.be/5D_A4IBWSv4?si=pgNinSqUyLBks6po',
  'https://youtu.be/MoHnffhBwqs?si=_FGX4ucMtOTcD2to',
  'https://youtu.be/Eg4LUvUjUWI?si=YqeuNfTh_iTuj-dP',
  'https://youtu.be/bcHoBDw4G10?si=auASu-G_c9NkS48Z',
  'https://youtu.be/CPAoMCo7tNw?si=2rEiXXCn6UcySUVZ',
  'https://youtu.be/EVDYmBrl02Q?si=ODB07HFZCtTtg4F4',
  'https://youtu.be/RBrdl0v_anc?si=cu3qNsVyUIIzZGvv',
  'https://youtu.be/50VNCymT-Cs?si=sEwBTlJCeuqL9LTD',
  'https://youtu.be/Cu5hhxP_prE?si=VRZVlVcLWqk8Dasg',
  'https://youtu.be/9gqAq6kq5Ek?si=Gro32XWDuPLWzyIv',
  'https://youtu.be/UprcpdwuwCg?si=O6_fwxx8TOkfjIXi',
  'https://youtu.be/Io0fBr1XBUA?si=SUp9MdCXlOU_Vf5s',
  'https://youtu.be/zABLecsR5UE?si=k3rryaA0P3O8JBhY',
  'https://youtu.be/n1NTv6Y4pxs?si=76WA3JI0TGILBHm7',
  'https://youtu.be/DeumyOzKqgI?si=Cok0dR7byK6pN682',
  'https://youtu.be/W-TE_Ys4iwM?si=RViOxRuaXxdz3pmm',
  'https://youtu.be/jtoncUzV6nA?si=yULSO1-MxnAVV13i',
  'https://www.youtube.com/watch?v=YykjpeuMNEk',
  'https://youtu.be/zgDbp5C74sU?si=R8Q5HZq2vzhGL57g',
  'https://www.youtube.com/watch?v=9bCp7j3nC30',
  'https://youtu.be/h_-JFUci0BM?si=SHiuHs1NdIjpN0WP',
  'https://youtu.be/S_0r3hYg78o?si=Be6GShy7mgRcl9Ha',
  'https://youtu.be/gHp-OjLOG5A?si=0abUDswbKz6rhQeX',
  'https://youtu.be/oIKuyj2GQtY',
  'https://youtu.be/kBqqlW6-99M?si=kXaaJTqhA4PaY6Gd',
  'https://youtu.be/DeFWClW7skQ?si=hkIGl-CTTw-FbnLz',
  'https://www.youtube.com/watch?v=vN0gaXS8dQE',
  'https://youtu.be/3WmdZOF5bKk?si=LcXY8Gohxxx4cZSA',
  'https://youtu.be/kkzEs0gdvZI?si=Z456wgKuJd0aE_PA',
  'https://youtu.be/ctmS5XX67Ek?si=NGZGPw0bcpfZciyi',
  'https://youtu.be/Jkj36B1YuDU?si=Yku5tRPe7avRNr2R',
  'https://youtu.be/OkxVxox--Io?si=AE4wj_c_uqTWGrbB',
  'https://youtu.be/Qau6mObfSGM?si=RsrcZ0VUCOHaEwE4',
  'https://youtu.be/jqyJ4xW2gb0?si=VgrA4JKMWkeWDIA5',
  'https://youtu.be/hNd5pILkpSw?si=qiwZxiuS0yeiuOPs',
  'https://www.youtube.com/watch?v=PPkJeWPP2nM',
  'https://youtu.be/omvW1cI-3xg?si=zHiFadZaUUpddcgu',
  'https://youtu.be/HTcL9WkB_wg?si=ILXw9EaPM4GJyx29',
  'https://youtu.be/WNeLUngb-Xg?si=V95nGOt0sMvhQG7c',
  'https://youtu.be/GB_S2qFh5lU?si=XDH6CdXhqJq-g321',
  'https://youtu.be/CBfSeqfeggI?si=A2RRcM0PSvAsvpHF',
  'https://youtu.be/t1-yL-xvklc?si=YZ1rS5hZtleOFOy1',
  'https://youtu.be/qzn_6bXdgeE?si=BgnimyD5Frnn_-o-',
  'https://youtu.be/pm3rDbXbZRI?si=7TxDuViBxhHGeZoU',
  'https://youtu.be/qzn_6bXdgeE?si=BgnimyD5Frnn_-o-',
  'https://youtu.be/pm3rDbXbZRI?si=7TxDuViBxhHGeZoU',
  'https://youtu.be/gzFmctgW0s8?si=JfvG_0Sj-IdPWoBY'
)

$startRow = 2
for ($i = 0; $i -lt $titles.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 2).Value = $titles[$i]
  $ws.Cells.Item($r, 3).Value = $links[$i]
}

